# Update cryptocurrency price / volume data (auto-refresh from GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.120.87"
$ws.Range("E2").Value = "  -4.54%  "

$ws.Range("D3").Value = "2.979.46"
$ws.Range("E3").Value = "  -1.46%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'559.32"
$ws.Range("E5").Value = "  -3.40%  "

$ws.Range("D6").Value = "'133.93"
$ws.Range("E6").Value = "  +4.66%  "

$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("E8").Value = "  +2.90%  "

$ws.Range("D9").Value = "2.973.65"
$ws.Range("E9").Value = "  -1.62%  "

$ws.Range("E10").Value = "  -3.33%  "

$ws.Range("D11").Value = "'4.87"
$ws.Range("E11").Value = "  -6.21%  "

$ws.Range("E13").Value = "  -0.42%  "

$ws.Range("D14").Value = "'33.07"
$ws.Range("E14").Value = "  +0.34%  "

$ws.Range("E15").Value = "  +0.82%  "

$ws.Range("D16").Value = "3.467.13"
$ws.Range("E16").Value = "  -1.50%  "

$ws.Range("D17").Value = "'6.89"
$ws.Range("E17").Value = "  +6.37%  "

$ws.Range("D18").Value = "2.971.59"
$ws.Range("E18").Value = "  -1.56%  "

$ws.Range("D19").Value = "57.992.35"
$ws.Range("E19").Value = "  -4.54%  "

$ws.Range("D20").Value = "'421.49"
$ws.Range("E20").Value = "  -3.35%  "

$ws.Range("D21").Value = "'13.28"
$ws.Range("E21").Value = "  +0.31%  "

$ws.Range("E22").Value = "  +2.90%  "

$ws.Range("E23").Value = "  -0.83%  "

$ws.Range("D24").Value = "'13.08"
$ws.Range("E24").Value = "  +1.16%  "

$ws.Range("D25").Value = "'79.73"
$ws.Range("E25").Value = "  -0.04%  "

$ws.Range("E26").Value = "  -0.14%  "

$ws.Range("E27").Value = "  +0.02%  "

$ws.Range("D28").Value = "'2.50"
$ws.Range("E28").Value = "  -2.63%  "

$ws.Range("D29").Value = "'7.62"
$ws.Range("E29").Value = "  +3.10%  "

$ws.Range("D30").Value = "'2.02"
$ws.Range("E30").Value = "  +4.75%  "

$ws.Range("D31").Value = "'25.38"
$ws.Range("E31").Value = "  -0.35%  "

$ws.Range("D32").Value = "'6.11"
$ws.Range("E32").Value = "  -2.53%  "

$ws.Range("E33").Value = "  +5.47%  "

$ws.Range("E34").Value = "  -1.93%  "

$ws.Range("D35").Value = "'5.68"
$ws.Range("E35").Value = "  +0.34%  "

$ws.Range("E36").Value = "  -1.69%  "

$ws.Range("D37").Value = "0.0₃0706"
$ws.Range("E37").Value = "  +4.26%  "

$ws.Range("D38").Value = "'48.78"
$ws.Range("E38").Value = "  -2.71%  "

$ws.Range("D39").Value = "'8.65"
$ws.Range("E39").Value = "  +1.01%  "

$ws.Range("D40").Value = "'2.60"
$ws.Range("E40").Value = "  +1.82%  "

$ws.Range("E41").Value = "  -3.31%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.108"
$ws.Range("E42").Value = "  -2.36%  "

$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "'379.62"
$ws.Range("E43").Value = "  -2.15%  "

$ws.Range("D44").Value = "2.692.55"
$ws.Range("E44").Value = "  +0.86%  "

$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("E46").Value = "  +2.44%  "

$ws.Range("D47").Value = "'121.90"
$ws.Range("E47").Value = "  +2.77%  "

$ws.Range("E48").Value = "  +2.54%  "

$ws.Range("E49").Value = "  -2.57%  "

$ws.Range("D50").Value = "'23.68"
$ws.Range("E50").Value = "  -1.94%  "

$ws.Range("E51").Value = "  -0.80%  "
